# feat: add 2022-Q1 data
#
# Before: sheets [2020-Q4, 总计]
# After:  sheets [2020-Q4, 2022-Q1, 总计]
#   - new "2022-Q1" sheet gets the per-fund holdings table (7 data columns)
#   - "总计" sheet gets a new summary row for 2022-Q1, above the existing 2020-Q4 row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet between "2020-Q4" and "总计".
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# Re-fetch "总计" - inserting a sheet shifts sheet positions, so any handle
# obtained before the Add() call now points at the wrong sheet.
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 2) Build the "2022-Q1" sheet.
#    Seed it from 总计's current A1:D2 (header style s=2 + indexed-row style
#    s=2 on column A) so the new sheet matches the workbook's existing look,
#    then stretch the header style across the extra columns E:H.
# ---------------------------------------------------------------------------
$total.Range("A1:D2").Copy($newSheet.Range("A1"))
$newSheet.Range("D1").Copy($newSheet.Range("E1:H1"))
$newSheet.Range("A2").Copy($newSheet.Range("A3"))

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - 银华-道琼斯88指数 A
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'180003"
$newSheet.Range("C2").Value = "银华-道琼斯88指数 A"
$newSheet.Range("D2").Value = "'14.84"
$newSheet.Range("E2").Value = "'83.54"
$newSheet.Range("F2").Value = "'2.90"
$newSheet.Range("G2").Value = "'0.4304"
$newSheet.Range("H2").Value = 9

# Row 3 - 银华沪港深增长股票
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'001703"
$newSheet.Range("C3").Value = "银华沪港深增长股票"
$newSheet.Range("D3").Value = "'3.02"
$newSheet.Range("E3").Value = "'87.52"
$newSheet.Range("F3").Value = "'2.82"
$newSheet.Range("G3").Value = "'0.0852"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------------
# 3) Update "总计": push the existing 2020-Q4 row down to row 3 and add the
#    new 2022-Q1 summary row at row 2.
# ---------------------------------------------------------------------------
$total.Range("A2:D2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.52

# Keep the original active sheet/selection ("2020-Q4" was active before this edit).
$wb.Worksheets.Item("2020-Q4").Activate() | Out-Null
$wb.Worksheets.Item("2020-Q4").Range("A1").Select() | Out-Null
